# The worksheet "type kinds" gains a new "Implemented meta interfaces"
# column. In the authored edit this is produced by moving the existing
# column F (which already held the "Implemented meta interfaces" data)
# to sit right before column D, shifting the old D/E columns one place
# to the right (D->E, E->F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Columns("F").Cut()
[void]$ws.Columns("D").Insert()

# Reflect the post-edit selection/scroll state: column E (the former
# column D) ends up selected, and the sheet view scrolls back to the top.
[void]$ws.Columns("E").Select()

# Rename the worksheet to reflect the broader topic of the table.
$ws.Name = "meta types"
